$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-12-20 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-21 Sunday", 2)

# Update the five rows of division problems in the single table.
# Row 1 (data row 1 of 5)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "86÷7=12, 2"
$t.Cell(1, 2).Range.Text = "83÷6=13, 5"
$t.Cell(1, 3).Range.Text = "71÷5=14, 1"
$t.Cell(1, 4).Range.Text = "32÷5=6, 2"
$t.Cell(1, 5).Range.Text = "56÷4=14, 0"

# Row 2 (data row 2 of 5 -> table row 5)
$t.Cell(5, 1).Range.Text = "61÷9=6, 7"
$t.Cell(5, 2).Range.Text = "94÷8=11, 6"
$t.Cell(5, 3).Range.Text = "63÷9=7, 0"
$t.Cell(5, 4).Range.Text = "80÷5=16, 0"
$t.Cell(5, 5).Range.Text = "17÷3=5, 2"

# Row 3 (data row 3 of 5 -> table row 9)
$t.Cell(9, 1).Range.Text = "25÷7=3, 4"
$t.Cell(9, 2).Range.Text = "34÷5=6, 4"
$t.Cell(9, 3).Range.Text = "53÷7=7, 4"
$t.Cell(9, 4).Range.Text = "93÷4=23, 1"
$t.Cell(9, 5).Range.Text = "76÷6=12, 4"

# Row 4 (data row 4 of 5 -> table row 13)
$t.Cell(13, 1).Range.Text = "75÷2=37, 1"
$t.Cell(13, 2).Range.Text = "70÷9=7, 7"
$t.Cell(13, 3).Range.Text = "26÷7=3, 5"
$t.Cell(13, 4).Range.Text = "88÷9=9, 7"
$t.Cell(13, 5).Range.Text = "99÷5=19, 4"

# Row 5 (data row 5 of 5 -> table row 17)
$t.Cell(17, 1).Range.Text = "91÷3=30, 1"
$t.Cell(17, 2).Range.Text = "40÷6=6, 4"
$t.Cell(17, 3).Range.Text = "19÷2=9, 1"
$t.Cell(17, 4).Range.Text = "89÷7=12, 5"
$t.Cell(17, 5).Range.Text = "57÷4=14, 1"
